$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Göstergelerdeki TFT projelerinden 2 sinin proje kodu güncellendi.
# Row 28: Model-01 TFT project, RNS Mcu -> RNS+Nuvo Mcu, code RNS -> ACM
$ws.Range("B28").Value = "DP-000-TF-ACM-CLI-H3B1-01"
$ws.Range("E28").Value = "RNS+Nuvo Mcu"

# Row 43: Model-02 TFT project, RNS Mcu -> RNS+Nuvo Mcu, code RNS -> ACM
$ws.Range("B43").Value = "DP-000-TF-ACM-CLI-H3B1-02"
$ws.Range("E43").Value = "RNS+Nuvo Mcu"
$ws.Range("K43").Value = "https://github.com/btk42/DP-000-TF-ACM-CLI-H3B1-02"

# Column E widened slightly to fit the new longer "RNS+Nuvo Mcu" text
$ws.Range("E1").ColumnWidth = 20.88671875

# Move the saved view/selection state to match the author's session
$ws.Application.ActiveWindow.ScrollRow = 22
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("M31").Select()

# Turn off iterative calculation (was enabled before)
$excel.Iteration = $false
